# Append the 2025-09-09 tracker entries (rows 12-16) to the progress history sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$goals = @(
    @("G2", "Workout"),
    @("G3", "Eat Healthy"),
    @("G4", "Read Book"),
    @("G5", "Investment Plan"),
    @("G6", "Spend 10 Hours without phone")
)

$startRow  = 12
$dateValue = 45909
$progress  = 0.9802960494069208
$percentage = 0
$change    = -0.01

for ($i = 0; $i -lt $goals.Count; $i++) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).Value = $goals[$i][0]
    $ws.Cells.Item($r, 2).Value = $goals[$i][1]

    $ws.Cells.Item($r, 3).Value = $dateValue
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 4).Value = $progress
    $ws.Cells.Item($r, 5).Value = $percentage
    $ws.Cells.Item($r, 6).Value = $change
}
